# Automatische test-sync: 2025-08-01 23:43:50
#
# Adds a new "Testmail #4" row to the Logs sheet, updates the matching
# category tally on the Dashboard sheet, extends the conditional
# formatting ranges on Logs to cover the new row, and extends the bar
# chart's category/value series references to include the new Dashboard
# row.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Logs: append row 9 -----------------------------------------------
$wsLogs.Range("A9").Value = "Wil je 100 stuks M5-bouten bestellen?"
$wsLogs.Range("B9").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C9").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$wsLogs.Range("D9").Value = "Bestelling / Levering"
$wsLogs.Range("E9").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$wsLogs.Range("F9").Value = "2025-08-01 23:43:02"
$wsLogs.Range("G9").Value = "Ja"
$wsLogs.Range("H9").Value = "Ja"
$wsLogs.Range("I9").Value = "Nee"
$wsLogs.Range("J9").Value = "Nee"

# --- Logs: extend conditional formatting ranges to include row 9 ------
$ranges = @("D2:D8", "G2:G8", "H2:H8", "I2:I8", "J2:J8")
$newRanges = @("D2:D9", "G2:G9", "H2:H9", "I2:I9", "J2:J9")

for ($i = 0; $i -lt $ranges.Count; $i++) {
    $oldRange = $ranges[$i]
    $newRange = $newRanges[$i]
    $fcs = $wsLogs.Range($oldRange).FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fc = $fcs.Item($j)
        $fc.ModifyAppliesToRange($wsLogs.Range($newRange))
    }
}

# --- Dashboard: append row 4 (tally for "Bestelling / Levering") ------
$wsDash.Range("A4").Value = "Bestelling / Levering"
$wsDash.Range("B4").Value = 1

# --- Chart: extend category/value series references to row 4 ---------
$co = $wsDash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
